$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Home_Score (J) and Away_Score (K) values for rows 10-19
$scores = @{
    10 = @(0, 1)
    11 = @(0, 1)
    12 = @(3, 1)
    13 = @(2, 1)
    14 = @(2, 2)
    15 = @(2, 0)
    16 = @(1, 1)
    17 = @(1, 1)
    18 = @(1, 1)
    19 = @(1, 0)
}

foreach ($row in $scores.Keys) {
    $values = $scores[$row]
    $ws.Cells.Item($row, 10).Value = $values[0]
    $ws.Cells.Item($row, 11).Value = $values[1]
}
